# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) is recomputed from upstream stats (std/mean of the
# underlying strike-price distribution) and rewritten here with the
# freshly-calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(6,6,7,4,8,4,8,2,2,4,3,3,1,4,2,5,2,6,5,5,3,3,13,4,9,4,4,12,1,3,7,7,4,3,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
